$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 395, shifting existing rows 395-418 down to 396-419.
$ws.Rows.Item(395).Insert()

# Populate the newly inserted row 395 with the new record.
$ws.Range("A395").Value2 = 6
$ws.Range("B395").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C395").Value2 = "Metropolitana"
$ws.Range("D395").Value2 = 44753
$ws.Range("E395").Value2 = 13
$ws.Range("F395").Value2 = 100112032
$ws.Range("G395").Value2 = "Zapallo italiano"
$ws.Range("H395").Value2 = "Sin especificar"
$ws.Range("I395").Value2 = "Primera"
$ws.Range("J395").Value2 = 430
$ws.Range("K395").Value2 = 9000
$ws.Range("L395").Value2 = 10000
$ws.Range("M395").Value2 = 9395
$ws.Range("N395").Value2 = "`$/caja 50 unidades"
$ws.Range("O395").Value2 = "Región de Arica y Parinacota"
$ws.Range("P395").Value2 = 188
$ws.Range("Q395").Value2 = 50
$ws.Range("R395").Value2 = "Hortaliza"
